$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 20330051920237
$ws.Range("B2").Value = "MARIANO"
$ws.Range("C2").Value = "ANTONIO"
$ws.Range("D2").Value = "JAQUELINE"
$ws.Range("E2").Value = "TOMA MUESTRAS BIOLÓGICAS"
$ws.Range("F2").Value = "2ALCM"
$ws.Range("G2").Value = 2
